# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Mon Nov 25 23:58:39 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so that numeric-looking strings
# (e.g. "3.20", "0.0000248", thousand-separated prices like "93.161.33")
# and padded percentages are preserved exactly as text, matching the
# original inlineStr cell contents instead of being coerced to numbers.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '93.161.33'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -4.88%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.419.09'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '234.81'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -7.13%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '637.71'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -3.40%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.42'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -1.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.394'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -8.40%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -7.36%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '3.418.69'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.75%  '
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -6.00%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '41.45'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -1.39%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.09'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.80%  '
$ws.Range('B15').NumberFormat = "@"
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').NumberFormat = "@"
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.062.24'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +1.94%  '
$ws.Range('B16').NumberFormat = "@"
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').NumberFormat = "@"
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '92.905.67'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -4.96%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000248'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -3.59%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '8.25'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -6.35%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.420.16'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.52%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.40'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -2.79%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.15'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +2.41%  '
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -9.65%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '491.76'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -4.42%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.20'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -5.70%  '
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -5.79%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.34'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -7.64%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '90.39'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -6.84%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.600.40'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +1.58%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '11.84'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -5.36%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '11.52'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.98%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.72'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +4.26%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.135'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -8.20%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -7.73%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '29.83'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +3.55%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.548'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -4.65%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '545.74'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +3.64%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -6.32%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '7.51'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -6.58%  '
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.906'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +5.26%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '23.97'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -1.97%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.70'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -3.77%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0405'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -10.56%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -4.47%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.52'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -3.05%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +2.94%  '
$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.18'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.52%  '
$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '52.72'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -2.45%  '
